# Edit script: apply "Burger King.xlsx" commit changes via Excel COM automation.
# Target changes (see task diff):
#  1. Sheet "Treinamento" (sheet1): add header cell B1 = "B" (new shared string).
#     Final selection/active sheet ends on Treinamento, cell B2 selected.
#  2. Sheet "Teste" (sheet2): change B1 from numeric 1 to text "B" (shared string,
#     same text as sheet1's new header); flip 66 label values in column B
#     (re-labelled training examples); append a new row 154 with only B154 = 1
#     (no A154 text). Final selection on that sheet is also B2 (but it is no
#     longer the active/tab-selected sheet).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Treinamento"
$ws2 = $wb.Worksheets.Item(2)   # "Teste"

# --- New header label used on both sheets --------------------------------
$ws1.Cells.Item(1, 2).Value = "B"
$ws2.Cells.Item(1, 2).Value = "B"

# --- Re-labelled rows on the "Teste" sheet (column B, 0/1 flips) ---------
$ws2.Cells.Item(5, 2).Value = 1
$ws2.Cells.Item(8, 2).Value = 0
$ws2.Cells.Item(9, 2).Value = 1
$ws2.Cells.Item(10, 2).Value = 0
$ws2.Cells.Item(11, 2).Value = 1
$ws2.Cells.Item(12, 2).Value = 0
$ws2.Cells.Item(13, 2).Value = 1
$ws2.Cells.Item(15, 2).Value = 0
$ws2.Cells.Item(17, 2).Value = 1
$ws2.Cells.Item(20, 2).Value = 0
$ws2.Cells.Item(23, 2).Value = 1
$ws2.Cells.Item(24, 2).Value = 0
$ws2.Cells.Item(25, 2).Value = 1
$ws2.Cells.Item(29, 2).Value = 0
$ws2.Cells.Item(30, 2).Value = 1
$ws2.Cells.Item(31, 2).Value = 0
$ws2.Cells.Item(32, 2).Value = 1
$ws2.Cells.Item(37, 2).Value = 0
$ws2.Cells.Item(40, 2).Value = 1
$ws2.Cells.Item(41, 2).Value = 0
$ws2.Cells.Item(44, 2).Value = 1
$ws2.Cells.Item(47, 2).Value = 0
$ws2.Cells.Item(48, 2).Value = 1
$ws2.Cells.Item(50, 2).Value = 0
$ws2.Cells.Item(51, 2).Value = 1
$ws2.Cells.Item(59, 2).Value = 0
$ws2.Cells.Item(62, 2).Value = 1
$ws2.Cells.Item(66, 2).Value = 0
$ws2.Cells.Item(69, 2).Value = 1
$ws2.Cells.Item(72, 2).Value = 0
$ws2.Cells.Item(76, 2).Value = 1
$ws2.Cells.Item(77, 2).Value = 0
$ws2.Cells.Item(78, 2).Value = 1
$ws2.Cells.Item(87, 2).Value = 0
$ws2.Cells.Item(88, 2).Value = 1
$ws2.Cells.Item(89, 2).Value = 0
$ws2.Cells.Item(91, 2).Value = 1
$ws2.Cells.Item(97, 2).Value = 0
$ws2.Cells.Item(100, 2).Value = 1
$ws2.Cells.Item(101, 2).Value = 0
$ws2.Cells.Item(102, 2).Value = 1
$ws2.Cells.Item(103, 2).Value = 0
$ws2.Cells.Item(108, 2).Value = 1
$ws2.Cells.Item(111, 2).Value = 0
$ws2.Cells.Item(114, 2).Value = 1
$ws2.Cells.Item(116, 2).Value = 0
$ws2.Cells.Item(119, 2).Value = 1
$ws2.Cells.Item(120, 2).Value = 0
$ws2.Cells.Item(127, 2).Value = 1
$ws2.Cells.Item(128, 2).Value = 0
$ws2.Cells.Item(130, 2).Value = 1
$ws2.Cells.Item(131, 2).Value = 0
$ws2.Cells.Item(132, 2).Value = 1
$ws2.Cells.Item(133, 2).Value = 0
$ws2.Cells.Item(134, 2).Value = 1
$ws2.Cells.Item(135, 2).Value = 0
$ws2.Cells.Item(136, 2).Value = 1
$ws2.Cells.Item(137, 2).Value = 0
$ws2.Cells.Item(138, 2).Value = 1
$ws2.Cells.Item(139, 2).Value = 0
$ws2.Cells.Item(141, 2).Value = 1
$ws2.Cells.Item(143, 2).Value = 0
$ws2.Cells.Item(144, 2).Value = 1
$ws2.Cells.Item(148, 2).Value = 0
$ws2.Cells.Item(150, 2).Value = 1
$ws2.Cells.Item(153, 2).Value = 0

# --- New trailing row: Teste!B154 = 1, no A154 text -----------------------
$ws2.Cells.Item(154, 2).Value = 1

# --- Final selections / active sheet -------------------------------------
# Teste ends up scrolled/selected at B2 but is no longer the selected tab.
$null = $ws2.Range("B2").Select()

# Treinamento becomes the selected tab, also parked at B2.
$null = $ws1.Activate()
$null = $ws1.Range("B2").Select()
